# Add five new AODVv2 issues (#58-#62) from the IETF issue tracker to both
# the "MANET-all" sheet (sheet1) and the "AODVv2-only" sheet (sheet2).

$wb = $excel.ActiveWorkbook

$manet  = $wb.Worksheets.Item(1)   # "MANET-all"
$aodvv2 = $wb.Worksheets.Item(2)   # "AODVv2-only"

# New ticket data: Ticket#, Summary, Component, Type, Owner, Status, Created
$newIssues = @(
    @{ Ticket = "#58"; Num = 58; Summary = 'Definitions of OrigNode and TargNode (Submitted for Justin Dean)'; Component = "aodvv2"; Type = "enhancement"; Owner = "draft-ietf-manet-aodvv2@tools.ietf.org"; Status = "new"; Created = 41991 },
    @{ Ticket = "#59"; Num = 59; Summary = 'Use of the term "invalid" (Submitted for Justin Dean)'; Component = "aodvv2"; Type = "enhancement"; Owner = "Charles Perkins"; Status = "new"; Created = 41991 },
    @{ Ticket = "#60"; Num = 60; Summary = "Should OrigNode be included in the message header? (Submitted for Justin Dean)"; Component = "aodvv2"; Type = "defect"; Owner = "Charles Perkins"; Status = "new"; Created = 41991 },
    @{ Ticket = "#61"; Num = 61; Summary = 'Difference between "broken" and "expired" (Submitted for Justin Dean)'; Component = "aodvv2"; Type = "defect"; Owner = "Charles Perkins"; Status = "new"; Created = 41991 },
    @{ Ticket = "#62"; Num = 62; Summary = 'Inconsistency surrounding the "timed" state (Submitted for Justin Dean)'; Component = "aodvv2"; Type = "enhancement"; Owner = "Charles Perkins"; Status = "new"; Created = 41991 }
)

# --- Sheet1 "MANET-all": columns A Ticket, B Summary, C Component, D Version,
#     E Milestone, F Type, G Owner, H Status, I Created. Starts at row 59.
$row = 59
foreach ($issue in $newIssues) {
    $manet.Range("A$row").Value = $issue.Ticket
    $manet.Range("B$row").Value = $issue.Summary
    $manet.Range("C$row").Value = $issue.Component
    $manet.Range("F$row").Value = $issue.Type
    $manet.Range("G$row").Value = $issue.Owner
    $manet.Range("H$row").Value = $issue.Status
    $manet.Range("I$row").Value = $issue.Created
    $manet.Range("I$row").NumberFormat = "m/d/yy"
    $row++
}

# --- Sheet2 "AODVv2-only": columns A Ticket (numeric), B Summary, C Component,
#     D Type, E Owner, F Status, G Created. Starts at row 52.
$row = 52
foreach ($issue in $newIssues) {
    $aodvv2.Range("A$row").Value = $issue.Num
    $aodvv2.Range("B$row").Value = $issue.Summary
    $aodvv2.Range("C$row").Value = $issue.Component
    $aodvv2.Range("D$row").Value = $issue.Type
    $aodvv2.Range("E$row").Value = $issue.Owner
    $aodvv2.Range("F$row").Value = $issue.Status
    $aodvv2.Range("G$row").Value = $issue.Created
    $aodvv2.Range("G$row").NumberFormat = "m/d/yy"
    $row++
}

# Restore the view/selection state seen after the edit: MANET-all scrolled
# down with the new rows selected, then AODVv2-only re-activated as the
# visible tab with its own new-row selection.
$manet.Activate()
[void]$manet.Range("A59:I63").Select()

$aodvv2.Activate()
[void]$aodvv2.Range("A57").Select()
